$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were added to the "Fruta, Vega Modelo de Temuco - Plátano"
# log, inserted right before the existing row 309 (shifting everything below it
# down by two rows, i.e. old row 309 -> new row 311, ..., old row 331 -> new row 333).
$ws.Rows("309:310").Insert()

# --- New row 309 ---
$r = 309
$ws.Cells.Item($r, 1).Value = 10
$ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value = "La Araucanía"
$ws.Cells.Item($r, 4).Value = 44461
$ws.Cells.Item($r, 5).Value = 9
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108006
$ws.Cells.Item($r, 10).Value = "Plátano"
$ws.Cells.Item($r, 11).Value = "Barraganete"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 80
$ws.Cells.Item($r, 14).Value = 24000
$ws.Cells.Item($r, 15).Value = 24000
$ws.Cells.Item($r, 16).Value = 24000
$ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item($r, 18).Value = "Ecuador"
$ws.Cells.Item($r, 19).Value = 1200
$ws.Cells.Item($r, 20).Value = 20

# --- New row 310 ---
$r = 310
$ws.Cells.Item($r, 1).Value = 10
$ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value = "La Araucanía"
$ws.Cells.Item($r, 4).Value = 44461
$ws.Cells.Item($r, 5).Value = 9
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108006
$ws.Cells.Item($r, 10).Value = "Plátano"
$ws.Cells.Item($r, 11).Value = "Sin especificar"
$ws.Cells.Item($r, 12).Value = "Pintón"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 18000
$ws.Cells.Item($r, 16).Value = 18000
$ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item($r, 18).Value = "Ecuador"
$ws.Cells.Item($r, 19).Value = 900
$ws.Cells.Item($r, 20).Value = 20
